$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# Append new program rows with their "Choose" values
$ws.Range("A8").Value = "RWTH_Aachen_DDS"
$ws.Range("B8").Value = "Yes"

$ws.Range("A9").Value = "RWTH_Aachen_TIME"
$ws.Range("B9").Value = "Yes"

$ws.Range("A10").Value = "Uni_Goettingen_Applied_CS"
$ws.Range("B10").Value = "Yes"

# Extend the data validation (dropdown list) over the new rows
$ws.Range("B1:B10").Validation.Delete()
$ws.Range("B1:B10").Validation.Add(3, 1, 1, '"Yes,No"')

# Move the active selection to the next empty row, like Excel would after data entry
$ws.Range("A11").Select() | Out-Null
